$wb = $excel.ActiveWorkbook

# Sheet 1: PayNowCC_27
$ws = $wb.Worksheets.Item("PayNowCC_27")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Sat Aug 30 00:15:44 IST 2025"

# Sheet 2: PayNowSCFCC_27
$ws = $wb.Worksheets.Item("PayNowSCFCC_27")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Sat Aug 30 00:18:21 IST 2025"

# Sheet 3: PayNowDCFCC_27
$ws = $wb.Worksheets.Item("PayNowDCFCC_27")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Sat Aug 30 00:12:10 IST 2025"

# Sheet 4: CCDeferredCC_27
$ws = $wb.Worksheets.Item("CCDeferredCC_27")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Fri Aug 29 23:59:23 IST 2025"

# Sheet 5: CMCAutopayCC_27
$ws = $wb.Worksheets.Item("CMCAutopayCC_27")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Sat Aug 30 00:08:24 IST 2025"
